$wb = $excel.ActiveWorkbook

# --- Typography sheet: set "Wildcard Characters" (column G) for the
# Default/Large/Small/Xlarge/XXLarge rows, and fill in the missing
# "Wildcard Ranges" (column I) for the Small row.
$wsTypo = $wb.Worksheets.Item("Typography")

$wsTypo.Range("G4").Value = '.",'
$wsTypo.Range("G5").Value = '.",'
$wsTypo.Range("G6").Value = '.",'
$wsTypo.Range("I6").Value = "0-9"
$wsTypo.Range("G7").Value = '.",'
$wsTypo.Range("G8").Value = '.",'

# --- Translation sheet: add two new status message rows (save/load
# error handling while reading/writing settings to non-volatile memory).
$wsTrans = $wb.Worksheets.Item("Translation")

$wsTrans.Range("B53").Value = "STATUSMSG_SETTINGS_LOAD_ERR"
$wsTrans.Range("C53").Value = "Default"
$wsTrans.Range("D53").Value = "Center"
$wsTrans.Range("E53").Value = "LTR"
$wsTrans.Range("F53").Value = "Settings could not be loaded from non-volatile memory. Manually set parameters."
$wsTrans.Range("G53").Value = "Nastavitve niso bile uspesno vnesene. Nastavi rocno,"

$wsTrans.Range("B54").Value = "STATUSMSG_SETTINGS_SAVE_ERR"
$wsTrans.Range("C54").Value = "Default"
$wsTrans.Range("D54").Value = "Center"
$wsTrans.Range("E54").Value = "LTR"
$wsTrans.Range("F54").Value = "Error writing to memory. Try again or reset the device."
$wsTrans.Range("G54").Value = "Pisanje v spomin je bilo neuspesno. Poskusi ponovno ali resetiraj napravo."
